$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New task rows to append below the existing data (row 14 already has "Bugs: Figure out css alignment")
# Row 14 gains a "Done" status in column D.
$ws.Range("D14").Value = "Done"

# New rows 15-20, skip row 21 (blank), then rows 22-23
$ws.Range("B15").Value = "Refactor components such that we have a gamestate component"
$ws.Range("D15").Value = "Done"

$ws.Range("B16").Value = "API Function to buy Card"
$ws.Range("B17").Value = "API Function to buy Resource"
$ws.Range("B18").Value = "API Function to buy Generators"
$ws.Range("B19").Value = "Add Map Component"
$ws.Range("B20").Value = "Add Domain classes for maps"

# Row 21 intentionally left blank

$ws.Range("B22").Value = "Need to cleanup angular services"
$ws.Range("B23").Value = "Implement the constant gamestate"

# Update selection to match the new active cell location
$ws.Range("B23").Select()
